$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.022.95"
$ws.Range("E2").Value = "  +3.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.228.32"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.48"
$ws.Range("E5").Value = "  +6.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "629.61"
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("E7").Value = "  +8.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.699"
$ws.Range("E8").Value = "  +5.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.227.11"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +7.53%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  +7.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.42"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.826.70"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.49"
$ws.Range("E16").Value = "  +4.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.756.19"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.235.87"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000233"
$ws.Range("E19").Value = "  +81.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.38"
$ws.Range("E20").Value = "  +15.19%  "
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "440.03"
$ws.Range("E22").Value = "  +6.52%  "
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.19"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.86"
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "80.52"
$ws.Range("E27").Value = "  +10.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.384.14"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.159"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.09"
$ws.Range("E32").Value = "  +37.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.59"
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "541.69"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.01"
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.48"
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.38"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.130"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "174.11"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.79"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.754"
$ws.Range("E48").Value = "  +10.13%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  +6.71%  "
